$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81, shifting existing rows 81:94 down to 82:95
$ws.Rows.Item(81).Insert()

# Match the date cell's number format to the one used elsewhere in column D
$ws.Cells.Item(81, 4).NumberFormat = $ws.Cells.Item(82, 4).NumberFormat

# Fill new row 81 data
$ws.Cells.Item(81, 1).Value = 8
$ws.Cells.Item(81, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(81, 3).Value = "Coquimbo"
$ws.Cells.Item(81, 4).Value = 44511
$ws.Cells.Item(81, 5).Value = 4
$ws.Cells.Item(81, 6).Value = 100112044
$ws.Cells.Item(81, 7).Value = "Perejil"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 3360
$ws.Cells.Item(81, 11).Value = 1300
$ws.Cells.Item(81, 12).Value = 1500
$ws.Cells.Item(81, 13).Value = 1400
$ws.Cells.Item(81, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(81, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(81, 16).Value = 933
$ws.Cells.Item(81, 17).Value = 1.5
$ws.Cells.Item(81, 18).Value = "Hortaliza"
